# Natmi LR-pairs output following Dr Hou advice:
# the Sema5a-Plxnb3 ligand-receptor table now includes "ECs" as a sending
# cluster too (previously it only appeared as a target cluster), giving a
# full 3x3 grid of ECs/FAPs/sCs x ECs/FAPs/sCs, and all the derived
# expression/specificity statistics are recomputed accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: ECs -> ECs
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Sema5a"
$ws.Range("C2").Value = "Plxnb3"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.4290636666666667
$ws.Range("H2").Value = 1.287191
$ws.Range("I2").Value = 0.0126431569814401
$ws.Range("J2").Value = 0.0126431569814401
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.432592
$ws.Range("N2").Value = 1.297776
$ws.Range("O2").Value = 0.1825161279706709
$ws.Range("P2").Value = 0.1825161279706708
$ws.Range("Q2").Value = 0.1856095096906667
$ws.Range("R2").Value = 1.670485587216
$ws.Range("S2").Value = 0.002307580057577802
$ws.Range("T2").Value = 0.002307580057577801

# Row 3: ECs -> FAPs
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Sema5a"
$ws.Range("C3").Value = "Plxnb3"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.4290636666666667
$ws.Range("H3").Value = 1.287191
$ws.Range("I3").Value = 0.0126431569814401
$ws.Range("J3").Value = 0.0126431569814401
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.3245703333333333
$ws.Range("N3").Value = 0.973711
$ws.Range("O3").Value = 0.1369403976359941
$ws.Range("P3").Value = 0.1369403976359941
$ws.Range("Q3").Value = 0.1392613373112222
$ws.Range("R3").Value = 1.253352035801
$ws.Range("S3").Value = 0.001731358944412702
$ws.Range("T3").Value = 0.001731358944412702

# Row 4: ECs -> sCs
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Sema5a"
$ws.Range("C4").Value = "Plxnb3"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.4290636666666667
$ws.Range("H4").Value = 1.287191
$ws.Range("I4").Value = 0.0126431569814401
$ws.Range("J4").Value = 0.0126431569814401
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 1.612995333333333
$ws.Range("N4").Value = 4.838986
$ws.Range("O4").Value = 0.6805434743933351
$ws.Range("P4").Value = 0.680543474393335
$ws.Range("Q4").Value = 0.6920776920362223
$ws.Range("R4").Value = 6.228699228326001
$ws.Range("S4").Value = 0.008604217979449593
$ws.Range("T4").Value = 0.008604217979449592

# Row 5: FAPs -> ECs
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Sema5a"
$ws.Range("C5").Value = "Plxnb3"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 26.436942
$ws.Range("H5").Value = 79.31082599999999
$ws.Range("I5").Value = 0.7790135445677298
$ws.Range("J5").Value = 0.7790135445677296
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.432592
$ws.Range("N5").Value = 1.297776
$ws.Range("O5").Value = 0.1825161279706709
$ws.Range("P5").Value = 0.1825161279706708
$ws.Range("Q5").Value = 11.436409613664
$ws.Range("R5").Value = 102.927686522976
$ws.Range("S5").Value = 0.1421825357912097
$ws.Range("T5").Value = 0.1421825357912097

# Row 6: FAPs -> FAPs
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Sema5a"
$ws.Range("C6").Value = "Plxnb3"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 26.436942
$ws.Range("H6").Value = 79.31082599999999
$ws.Range("I6").Value = 0.7790135445677298
$ws.Range("J6").Value = 0.7790135445677296
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 0.3245703333333333
$ws.Range("N6").Value = 0.973711
$ws.Range("O6").Value = 0.1369403976359941
$ws.Range("P6").Value = 0.1369403976359941
$ws.Range("Q6").Value = 8.580647077254
$ws.Range("R6").Value = 77.225823695286
$ws.Range("S6").Value = 0.1066784245569301
$ws.Range("T6").Value = 0.1066784245569301

# Row 7: FAPs -> sCs
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Sema5a"
$ws.Range("C7").Value = "Plxnb3"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 26.436942
$ws.Range("H7").Value = 79.31082599999999
$ws.Range("I7").Value = 0.7790135445677298
$ws.Range("J7").Value = 0.7790135445677296
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 1.612995333333333
$ws.Range("N7").Value = 4.838986
$ws.Range("O7").Value = 0.6805434743933351
$ws.Range("P7").Value = 0.680543474393335
$ws.Range("Q7").Value = 42.642664073604
$ws.Range("R7").Value = 383.783976662436
$ws.Range("S7").Value = 0.53015258421959
$ws.Range("T7").Value = 0.5301525842195898

# Row 8 (new): sCs -> ECs
$ws.Range("A8").Value = "sCs"
$ws.Range("B8").Value = "Sema5a"
$ws.Range("C8").Value = "Plxnb3"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 7.070428666666666
$ws.Range("H8").Value = 21.211286
$ws.Range("I8").Value = 0.2083432984508302
$ws.Range("J8").Value = 0.2083432984508302
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 0.432592
$ws.Range("N8").Value = 1.297776
$ws.Range("O8").Value = 0.1825161279706709
$ws.Range("P8").Value = 0.1825161279706708
$ws.Range("Q8").Value = 3.058610877770667
$ws.Range("R8").Value = 27.527497899936
$ws.Range("S8").Value = 0.03802601212188339
$ws.Range("T8").Value = 0.03802601212188339

# Row 9 (new): sCs -> FAPs
$ws.Range("A9").Value = "sCs"
$ws.Range("B9").Value = "Sema5a"
$ws.Range("C9").Value = "Plxnb3"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 7.070428666666666
$ws.Range("H9").Value = 21.211286
$ws.Range("I9").Value = 0.2083432984508302
$ws.Range("J9").Value = 0.2083432984508302
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 0.3245703333333333
$ws.Range("N9").Value = 0.973711
$ws.Range("O9").Value = 0.1369403976359941
$ws.Range("P9").Value = 0.1369403976359941
$ws.Range("Q9").Value = 2.294851389149555
$ws.Range("R9").Value = 20.653662502346
$ws.Range("S9").Value = 0.02853061413465128
$ws.Range("T9").Value = 0.02853061413465128

# Row 10 (new): sCs -> sCs
$ws.Range("A10").Value = "sCs"
$ws.Range("B10").Value = "Sema5a"
$ws.Range("C10").Value = "Plxnb3"
$ws.Range("D10").Value = "sCs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 7.070428666666666
$ws.Range("H10").Value = 21.211286
$ws.Range("I10").Value = 0.2083432984508302
$ws.Range("J10").Value = 0.2083432984508302
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 1.612995333333333
$ws.Range("N10").Value = 4.838986
$ws.Range("O10").Value = 0.6805434743933351
$ws.Range("P10").Value = 0.680543474393335
$ws.Range("Q10").Value = 11.40456844399955
$ws.Range("R10").Value = 102.641115995996
$ws.Range("S10").Value = 0.1417866721942955
$ws.Range("T10").Value = 0.1417866721942955

